# Add the "InsuranceExpiry" field to the vehicle upload template.
# It is inserted as a new column right after "InsuranceNumber" (column K),
# pushing DriverName..NormalSpeed one column to the right (L..R -> M..S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L; this shifts existing L:R -> M:S and copies
# the formatting of the column to the left (K, InsuranceNumber) onto it.
$ws.Columns("L").Insert()

# Give the new header cell its label.
$ws.Range("L1").Value = "InsuranceExpiry"

# Match the (auto-fit) column width Excel computed for the new header.
$ws.Columns("L").ColumnWidth = 14.5
